# Rename the header row: strip umlauts/capitalisation (old -> new)
#   Jahr            -> jahr
#   Umsatzerlöse    -> umsatzerloese
#   Verkehrserlöse  -> verkehrserloese
#   EBIT            -> ebit
#   Konzernergebnis -> konzernergebnis
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "jahr"
$ws.Range("B1").Value = "umsatzerloese"
$ws.Range("C1").Value = "verkehrserloese"
$ws.Range("D1").Value = "ebit"
$ws.Range("E1").Value = "konzernergebnis"

# Leave the selection where the author last left it when saving.
[void]$ws.Range("G10").Select()
